# Auto-refresh of resum_diari_meteocat.xlsx: data extraction timestamps (col E)
# plus a handful of re-sampled sensor readings (cols H, J, L, N, O) for the
# "2026-02-27 18:50" run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = '2026-02-27 18:48:42'

# Row 3
$ws.Range("E3").Value = '2026-02-27 18:48:45'

# Row 4
$ws.Range("E4").Value = '2026-02-27 18:48:48'

# Row 5
$ws.Range("E5").Value = '2026-02-27 18:48:50'
$ws.Range("H5").Value = '''38%'

# Row 6
$ws.Range("E6").Value = '2026-02-27 18:48:53'
$ws.Range("J6").Value = '1024.4 hPa'

# Row 7
$ws.Range("E7").Value = '2026-02-27 18:48:56'
$ws.Range("J7").Value = '1024.9 hPa'

# Row 8
$ws.Range("E8").Value = '2026-02-27 18:48:58'
$ws.Range("H8").Value = '''61%'
$ws.Range("N8").Value = '8.8 °C 18:25 TU'
$ws.Range("O8").Value = '12.2 °C'

# Row 9
$ws.Range("E9").Value = '2026-02-27 18:49:00'

# Row 10
$ws.Range("E10").Value = '2026-02-27 18:49:03'

# Row 11
$ws.Range("E11").Value = '2026-02-27 18:49:06'
$ws.Range("O11").Value = '8.9 °C'

# Row 12
$ws.Range("E12").Value = '2026-02-27 18:49:09'
$ws.Range("H12").Value = '''95%'

# Row 13
$ws.Range("E13").Value = '2026-02-27 18:49:11'
$ws.Range("H13").Value = '''66%'
$ws.Range("J13").Value = '1025.7 hPa'
$ws.Range("O13").Value = '6.6 °C'

# Row 14
$ws.Range("E14").Value = '2026-02-27 18:49:14'

# Row 15
$ws.Range("E15").Value = '2026-02-27 18:49:17'
$ws.Range("O15").Value = '10.9 °C'

# Row 16
$ws.Range("E16").Value = '2026-02-27 18:49:19'
$ws.Range("O16").Value = '2.9 °C'

# Row 17
$ws.Range("E17").Value = '2026-02-27 18:49:22'
$ws.Range("L17").Value = '46.4 km/h - 244º 18:20 TU'

# Row 18
$ws.Range("E18").Value = '2026-02-27 18:49:25'
$ws.Range("J18").Value = '1024.6 hPa'

# Row 19
$ws.Range("E19").Value = '2026-02-27 18:49:28'
$ws.Range("O19").Value = '10.6 °C'

# Row 20
$ws.Range("E20").Value = '2026-02-27 18:49:30'

# Row 21
$ws.Range("E21").Value = '2026-02-27 18:49:33'
$ws.Range("J21").Value = '1024.4 hPa'
$ws.Range("O21").Value = '9.8 °C'

# Row 22
$ws.Range("E22").Value = '2026-02-27 18:49:36'

# Row 23
$ws.Range("E23").Value = '2026-02-27 18:49:39'
$ws.Range("O23").Value = '3.9 °C'

# Row 24
$ws.Range("E24").Value = '2026-02-27 18:49:41'
$ws.Range("J24").Value = '1023.7 hPa'

# Row 25
$ws.Range("E25").Value = '2026-02-27 18:49:44'
$ws.Range("H25").Value = '''33%'

# Row 26
$ws.Range("E26").Value = '2026-02-27 18:49:47'
$ws.Range("H26").Value = '''43%'
$ws.Range("O26").Value = '10.6 °C'

# Row 27
$ws.Range("E27").Value = '2026-02-27 18:49:49'
$ws.Range("H27").Value = '''41%'

# Row 28
$ws.Range("E28").Value = '2026-02-27 18:49:52'
$ws.Range("J28").Value = '1024.7 hPa'
$ws.Range("O28").Value = '8.0 °C'

# Row 29
$ws.Range("E29").Value = '2026-02-27 18:49:55'

# Row 30
$ws.Range("E30").Value = '2026-02-27 18:49:57'

# Row 31
$ws.Range("E31").Value = '2026-02-27 18:50:00'

# Row 32
$ws.Range("E32").Value = '2026-02-27 18:50:03'

# Row 33
$ws.Range("E33").Value = '2026-02-27 18:50:05'
$ws.Range("J33").Value = '1023.9 hPa'

# Row 34
$ws.Range("E34").Value = '2026-02-27 18:50:08'
$ws.Range("H34").Value = '''46%'

# Row 35
$ws.Range("E35").Value = '2026-02-27 18:50:11'

# Row 36
$ws.Range("E36").Value = '2026-02-27 18:50:14'
$ws.Range("J36").Value = '1024.9 hPa'

# Row 37
$ws.Range("E37").Value = '2026-02-27 18:50:16'
$ws.Range("H37").Value = '''69%'
$ws.Range("J37").Value = '1025.0 hPa'
$ws.Range("O37").Value = '8.3 °C'

# Row 38
$ws.Range("E38").Value = '2026-02-27 18:50:19'

# Row 39
$ws.Range("E39").Value = '2026-02-27 18:50:21'

# Row 40
$ws.Range("E40").Value = '2026-02-27 18:50:24'
$ws.Range("J40").Value = '1024.9 hPa'
$ws.Range("O40").Value = '9.3 °C'

# Row 41
$ws.Range("E41").Value = '2026-02-27 18:50:27'

# Row 42
$ws.Range("E42").Value = '2026-02-27 18:50:30'

# Row 43
$ws.Range("E43").Value = '2026-02-27 18:50:32'
$ws.Range("H43").Value = '''74%'
$ws.Range("O43").Value = '9.5 °C'

# Row 44
$ws.Range("E44").Value = '2026-02-27 18:50:35'
$ws.Range("H44").Value = '''58%'

# Row 45
$ws.Range("E45").Value = '2026-02-27 18:50:38'
$ws.Range("J45").Value = '1021.8 hPa'

# Row 46
$ws.Range("E46").Value = '2026-02-27 18:50:41'
$ws.Range("J46").Value = '1024.3 hPa'
